$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("to do")
$ws.Activate()
$n = $wb.Names.Item(1)
$n.RefersTo = "='to do'!`$A`$1:`$H`$32"
Write-Output $n.RefersTo
